$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: make room for the "failure rate"/"Unit" columns that are
# moving over from Sheet2 by inserting two new columns at C:D. This shifts
# the existing Class/Division (and the extra numeric columns) two slots to
# the right (C/D -> E/F, E/F -> G/H) and grows the dimension to H.
[void]$ws1.Range("C1:D1").EntireColumn.Insert()

# --- Move the "failure rate" (C) / "Unit" (D) columns from Sheet2 into the
# newly freed C:D columns on Sheet1, preserving formulas/values and types.
for ($r = 2; $r -le 12; $r++) {
    $srcC = $ws2.Cells.Item($r, 3)
    $dstC = $ws1.Cells.Item($r, 3)
    if ($srcC.Formula -ne "") {
        $dstC.Formula = $srcC.Formula
    }

    $srcD = $ws2.Cells.Item($r, 4)
    $dstD = $ws1.Cells.Item($r, 4)
    if ($srcD.Formula -ne "") {
        $dstD.Formula = $srcD.Formula
    }
}

# --- Sheet2 no longer needs the old C:D data; it now only carries column A.
[void]$ws2.Range("C2:D12").ClearContents()

# --- Update the active sheet/selection bookkeeping to match: Sheet1 becomes
# the selected tab (with a new selection), Sheet2 loses tabSelected and gets
# a new selection pointing at the (now moved/cleared) C2:D12 block. Select
# on Sheet2 first (it's currently active) so the final Select() below is the
# one that sticks as "active".
[void]$ws2.Range("C2:D12").Select()

[void]$ws1.Activate()
[void]$ws1.Range("G17").Select()
